# The records that occupied rows 2-5 are cyclically rotated:
#   old row 5 -> new row 2
#   old row 2 -> new row 3
#   old row 3 -> new row 4
#   old row 4 -> new row 5
# (equivalent to lifting the row-5 record out and re-inserting it at row 2,
#  pushing the old rows 2-4 down by one). No other rows are affected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastCol = "AY"

# Snapshot the four rows first (must read everything before any writes).
$row2 = $ws.Range("A2:" + $lastCol + "2").Value2
$row3 = $ws.Range("A3:" + $lastCol + "3").Value2
$row4 = $ws.Range("A4:" + $lastCol + "4").Value2
$row5 = $ws.Range("A5:" + $lastCol + "5").Value2

$rotated = @{ 2 = $row5; 3 = $row2; 4 = $row3; 5 = $row4 }

# Columns Y and AA hold plain-text dates ("YYYY-MM-DD"). A bulk array write
# through .Value would make Excel auto-convert those strings into real date
# serials (and stamp a date NumberFormat on the cell). Blank those two
# columns out of the bulk payload and restore them afterwards as literal
# text so the cells stay General/number-format-free, matching the source.
$colY = 25
$colAA = 27

foreach ($r in 2..5) {
    $data = $rotated[$r]
    $yText = $data[1, $colY]
    $aaText = $data[1, $colAA]
    $data[1, $colY] = $null
    $data[1, $colAA] = $null

    $ws.Range("A" + $r + ":" + $lastCol + $r).Value = $data

    $yCell = $ws.Range("Y" + $r)
    if ($yText -eq $null) {
        $yCell.Value = $null
    } else {
        $yCell.Value = "'" + $yText
        $yCell.ClearFormats()
    }

    $aaCell = $ws.Range("AA" + $r)
    if ($aaText -eq $null) {
        $aaCell.Value = $null
    } else {
        $aaCell.Value = "'" + $aaText
        $aaCell.ClearFormats()
    }
}
